$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.657.60'
$ws.Range('E2').Value = '  +5.73%  '
$ws.Range('D3').Value = '3.634.14'
$ws.Range('E3').Value = '  +5.72%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.39'
$ws.Range('E5').Value = '  +2.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '194.96'
$ws.Range('E6').Value = '  +3.44%  '
$ws.Range('E7').Value = '  +2.65%  '
$ws.Range('D8').Value = '3.626.87'
$ws.Range('E8').Value = '  +5.82%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +8.07%  '
$ws.Range('E11').Value = '  +5.52%  '
$ws.Range('E12').Value = '  +1.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000314'
$ws.Range('E13').Value = '  +13.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.92'
$ws.Range('E14').Value = '  +5.27%  '
$ws.Range('D15').Value = '4.216.30'
$ws.Range('E15').Value = '  +5.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.53'
$ws.Range('E16').Value = '  +8.82%  '
$ws.Range('D17').Value = '3.635.49'
$ws.Range('E17').Value = '  +6.12%  '
$ws.Range('D18').Value = '70.667.15'
$ws.Range('E18').Value = '  +5.78%  '
$ws.Range('E19').Value = '  +5.80%  '
$ws.Range('E20').Value = '  +2.76%  '
$ws.Range('E21').Value = '  +4.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '488.28'
$ws.Range('E22').Value = '  +2.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.32'
$ws.Range('E23').Value = '  +14.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.25'
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.49'
$ws.Range('E25').Value = '  +3.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '91.30'
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('E27').Value = '  +6.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.41'
$ws.Range('E28').Value = '  +4.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.56'
$ws.Range('E29').Value = '  +6.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.96'
$ws.Range('E30').Value = '  +8.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.82'
$ws.Range('E31').Value = '  +5.62%  '
$ws.Range('E32').Value = '  +9.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.29'
$ws.Range('E33').Value = '  +4.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '66.26'
$ws.Range('E34').Value = '  +2.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '610.16'
$ws.Range('E35').Value = '  +2.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '40.53'
$ws.Range('E36').Value = '  +9.76%  '
$ws.Range('D37').Value = '0.0₃0840'
$ws.Range('E37').Value = '  +12.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.410'
$ws.Range('E38').Value = '  +5.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.148'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +3.03%  '
$ws.Range('D42').Value = '3.317.57'
$ws.Range('E42').Value = '  +4.07%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.17'
$ws.Range('E43').Value = '  +9.40%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.16'
$ws.Range('E44').Value = '  +17.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.86'
$ws.Range('E45').Value = '  +10.48%  '
$ws.Range('E46').Value = '  +6.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.66'
$ws.Range('E47').Value = '  +12.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.34'
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('E49').Value = '  +3.33%  '
$ws.Range('E50').Value = '  +1.53%  '
$ws.Range('E51').Value = '  -0.05%  '
